# Refresh the cryptos price/volume table with the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every "Price" (column D) cell in this sheet is stored as text in the
# source data (e.g. "232.63", "1.20") even though it looks numeric, so
# trailing/grouping-style digits survive. A plain .Value assignment lets
# Excel's automatic type detection silently coerce values like "14.90" or
# "6.00" back into real numbers (14.9, 6), so force each changed Price cell
# to text first (looping, since a multi-area Range here only applies
# NumberFormat to its first area), write the values, then restore the
# default "Normal" style so no visual formatting change is introduced.
$priceRefs = @("D2", "D3", "D5", "D6", "D8", "D10", "D12", "D13", "D14", "D15", "D17", "D18", "D19", "D20", "D21", "D22", "D24", "D25", "D26", "D27", "D28", "D29", "D32", "D33", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D49", "D50", "D51")
foreach ($ref in $priceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "37.471.07"
$ws.Range("E2").Value = "  -0.68%  "
$ws.Range("D3").Value = "2.073.02"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "232.19"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "57.51"
$ws.Range("E8").Value = "  -1.36%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Value = "14.90"
$ws.Range("E12").Value = "  +1.37%  "
$ws.Range("D13").Value = "2.382.36"
$ws.Range("E13").Value = "  +0.08%  "
$ws.Range("D14").Value = "20.93"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "0.767"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "2.086.44"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").Value = "37.398.48"
$ws.Range("E18").Value = "  -0.72%  "
$ws.Range("D19").Value = "70.40"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "6.00"
$ws.Range("E20").Value = "  -2.74%  "
$ws.Range("D21").Value = "0.0₃0828"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").Value = "228.17"
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "2.34"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.37"
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").Value = "9.63"
$ws.Range("E26").Value = "  +6.93%  "
$ws.Range("D27").Value = "169.98"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("D28").Value = "0.132"
$ws.Range("E28").Value = "  -3.99%  "
$ws.Range("D29").Value = "19.50"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("E30").Value = "  -1.28%  "
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").Value = "4.62"
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").Value = "0.0631"
$ws.Range("E33").Value = "  +0.17%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -2.32%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "5.27"
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "0.0230"
$ws.Range("E40").Value = "  +7.40%  "
$ws.Range("D41").Value = "99.80"
$ws.Range("E41").Value = "  -1.26%  "
$ws.Range("D42").Value = "2.91"
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").Value = "0.0954"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +4.08%  "
$ws.Range("D45").Value = "1.461.65"
$ws.Range("D46").Value = "16.64"
$ws.Range("E46").Value = "  +1.82%  "
$ws.Range("E47").Value = "  -1.20%  "
$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D48").Value = "3.95"
$ws.Range("E48").Value = "  -5.79%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "7.22"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("D50").Value = "2.95"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").Value = "2.266.60"
$ws.Range("E51").Value = "  +0.05%  "

# Restore default styling on the cells we force-formatted as text above.
foreach ($ref in $priceRefs) {
    $ws.Range($ref).Style = "Normal"
}
